# Auto-generated edits applying the diff to Sheets/Ramuh_Profits.xlsx
# Updates computed market-board profit columns (H,I,J,K,L,M,N) per leve row,
# per the scheduled-runner data refresh described in the commit message.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 80006
$ws.Range("J13").Value = 80006
$ws.Range("L13").Value = 80006
$ws.Range("N13").Value = -80344
$ws.Range("H33").Value = 165.07143
$ws.Range("I33").Value = 148.41667
$ws.Range("J33").Value = 265
$ws.Range("K33").Value = 148.41667
$ws.Range("L33").Value = 265
$ws.Range("M33").Value = 80.58332999999999
$ws.Range("N33").Value = -723

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 23968.666
$ws.Range("I3").Value = 950
$ws.Range("K3").Value = 950
$ws.Range("M3").Value = -835
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H12").Value = 1000
$ws.Range("J12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("N12").Value = -1346
$ws.Range("H13").Value = 2500588.2
$ws.Range("I13").Value = 10000003
$ws.Range("J13").Value = 783.3333
$ws.Range("K13").Value = 10000003
$ws.Range("L13").Value = 783.3333
$ws.Range("M13").Value = -9999859
$ws.Range("N13").Value = -1071.3333
$ws.Range("H14").Value = 49383028
$ws.Range("I14").Value = 49383028
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 49383028
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -49382853
$ws.Range("N14").ClearContents()
$ws.Range("H22").Value = 4492.857
$ws.Range("I22").Value = 1012.5
$ws.Range("J22").Value = 9133.333000000001
$ws.Range("K22").Value = 1012.5
$ws.Range("L22").Value = 9133.333000000001
$ws.Range("M22").Value = -713.5
$ws.Range("N22").Value = -9731.333000000001
$ws.Range("H25").Value = 375
$ws.Range("I25").Value = 375
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 375
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 27
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1160578.6
$ws.Range("I5").Value = 2750500.5
$ws.Range("J5").Value = 4271.8184
$ws.Range("K5").Value = 2750500.5
$ws.Range("L5").Value = 4271.8184
$ws.Range("M5").Value = -2750387.5
$ws.Range("N5").Value = -4497.8184
$ws.Range("H7").Value = 214571.42
$ws.Range("I7").Value = 214571.42
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 214571.42
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -214458.42
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 168
$ws.Range("I8").Value = 168
$ws.Range("K8").Value = 168
$ws.Range("M8").Value = -28
$ws.Range("H10").Value = 624.5
$ws.Range("I10").Value = 499.33334
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 499.33334
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -359.33334
$ws.Range("N10").Value = -1280
$ws.Range("H11").Value = 135.42857
$ws.Range("I11").Value = 135.42857
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 135.42857
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 4.571429999999992
$ws.Range("N11").ClearContents()
$ws.Range("H12").Value = 975.3333
$ws.Range("I12").Value = 970.4
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 970.4
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -802.4
$ws.Range("N12").Value = -1336
$ws.Range("H24").Value = 1325
$ws.Range("I24").Value = 1000
$ws.Range("J24").Value = 1650
$ws.Range("K24").Value = 1000
$ws.Range("L24").Value = 1650
$ws.Range("M24").Value = -765
$ws.Range("N24").Value = -2120
$ws.Range("H25").Value = 5830.4
$ws.Range("I25").Value = 801.3333
$ws.Range("J25").Value = 7985.7144
$ws.Range("K25").Value = 801.3333
$ws.Range("L25").Value = 7985.7144
$ws.Range("M25").Value = -566.3333
$ws.Range("N25").Value = -8455.714400000001
$ws.Range("H36").Value = 12885.25
$ws.Range("I36").Value = 1250
$ws.Range("J36").Value = 24520.5
$ws.Range("K36").Value = 1250
$ws.Range("L36").Value = 24520.5
$ws.Range("M36").Value = -716
$ws.Range("N36").Value = -25588.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 625
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 733.3333
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 733.3333
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -959.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 67.3
$ws.Range("I6").Value = 67.3
$ws.Range("K6").Value = 201.9
$ws.Range("M6").Value = -88.89999999999998
$ws.Range("H132").Value = 1046.8889
$ws.Range("I132").Value = 859.75
$ws.Range("J132").Value = 1319.091
$ws.Range("K132").Value = 7737.75
$ws.Range("L132").Value = 11871.819
$ws.Range("M132").Value = -5207.75
$ws.Range("N132").Value = -16931.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 740.75
$ws.Range("I3").Value = 321
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 321
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -205
$ws.Range("N3").Value = -2232
$ws.Range("H9").Value = 4500
$ws.Range("J9").Value = 4500
$ws.Range("L9").Value = 4500
$ws.Range("N9").Value = -4840
$ws.Range("H10").Value = 1002000
$ws.Range("I10").Value = 2000000
$ws.Range("K10").Value = 2000000
$ws.Range("M10").Value = -1999831
$ws.Range("H11").Value = 43752.5
$ws.Range("I11").Value = 21251.5
$ws.Range("J11").Value = 55003
$ws.Range("K11").Value = 21251.5
$ws.Range("L11").Value = 55003
$ws.Range("M11").Value = -21112.5
$ws.Range("N11").Value = -55281
$ws.Range("H13").Value = 371
$ws.Range("I13").Value = 213.75
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 213.75
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = -74.75
$ws.Range("N13").Value = -1278
$ws.Range("H23").Value = 500
$ws.Range("J23").Value = 500
$ws.Range("L23").Value = 500
$ws.Range("N23").Value = -946

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 60484
$ws.Range("I10").Value = 50964
$ws.Range("J10").Value = 70004
$ws.Range("K10").Value = 50964
$ws.Range("L10").Value = 70004
$ws.Range("M10").Value = -50824
$ws.Range("N10").Value = -70284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 500250
$ws.Range("I3").Value = 1000000
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 1000000
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = -999886
$ws.Range("N3").Value = -728
$ws.Range("H4").Value = 287811.5
$ws.Range("J4").Value = 500150
$ws.Range("L4").Value = 500150
$ws.Range("N4").Value = -500376
$ws.Range("H5").Value = 1117777.8
$ws.Range("I5").Value = 10000000
$ws.Range("J5").Value = 7500
$ws.Range("K5").Value = 10000000
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = -9999888
$ws.Range("N5").Value = -7724
$ws.Range("H6").Value = 940
$ws.Range("J6").Value = 940
$ws.Range("L6").Value = 940
$ws.Range("N6").Value = -1170
$ws.Range("H8").Value = 3000
$ws.Range("I8").Value = 3000
$ws.Range("K8").Value = 3000
$ws.Range("M8").Value = -2860
$ws.Range("H10").Value = 800
$ws.Range("J10").Value = 800
$ws.Range("L10").Value = 800
$ws.Range("N10").Value = -1138
$ws.Range("H11").Value = 1000
$ws.Range("J11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("N11").Value = -1284
$ws.Range("H31").Value = 5000
$ws.Range("J31").Value = 5000
$ws.Range("L31").Value = 5000
$ws.Range("N31").Value = -5696
